$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the species-record data (columns A, B, E, F, G, H, Q, R)
# among rows 2, 4 and 5, while leaving all other columns (C, D, P, S, T, U,
# V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) untouched, since they
# describe the same location/date/observer for every row.
#
# Rotation (old -> new):
#   row 2 <- old row 5
#   row 4 <- old row 2
#   row 5 <- old row 4

# Capture the original values before overwriting anything.
# NOTE: use .Value2 (not .Value) - in this runtime, reading/round-tripping
# .Value through a variable can yield the property descriptor instead of
# the actual cell content, while .Value2 reliably returns the raw value.
$row2 = @{
    A = $ws.Cells.Item(2, 1).Value2
    B = $ws.Cells.Item(2, 2).Value2
    E = $ws.Cells.Item(2, 5).Value2
    F = $ws.Cells.Item(2, 6).Value2
    G = $ws.Cells.Item(2, 7).Value2
    H = $ws.Cells.Item(2, 8).Value2
    Q = $ws.Cells.Item(2, 17).Value2
    R = $ws.Cells.Item(2, 18).Value2
}

$row4 = @{
    A = $ws.Cells.Item(4, 1).Value2
    B = $ws.Cells.Item(4, 2).Value2
    E = $ws.Cells.Item(4, 5).Value2
    F = $ws.Cells.Item(4, 6).Value2
    G = $ws.Cells.Item(4, 7).Value2
    H = $ws.Cells.Item(4, 8).Value2
    Q = $ws.Cells.Item(4, 17).Value2
    R = $ws.Cells.Item(4, 18).Value2
}

$row5 = @{
    A = $ws.Cells.Item(5, 1).Value2
    B = $ws.Cells.Item(5, 2).Value2
    E = $ws.Cells.Item(5, 5).Value2
    F = $ws.Cells.Item(5, 6).Value2
    G = $ws.Cells.Item(5, 7).Value2
    H = $ws.Cells.Item(5, 8).Value2
    Q = $ws.Cells.Item(5, 17).Value2
    R = $ws.Cells.Item(5, 18).Value2
}

# Write old row 5 values into row 2.
$ws.Cells.Item(2, 1).Value2 = $row5.A
$ws.Cells.Item(2, 2).Value2 = $row5.B
$ws.Cells.Item(2, 5).Value2 = $row5.E
$ws.Cells.Item(2, 6).Value2 = $row5.F
$ws.Cells.Item(2, 7).Value2 = $row5.G
$ws.Cells.Item(2, 8).Value2 = $row5.H
$ws.Cells.Item(2, 17).Value2 = $row5.Q
$ws.Cells.Item(2, 18).Value2 = $row5.R

# Write old row 2 values into row 4.
$ws.Cells.Item(4, 1).Value2 = $row2.A
$ws.Cells.Item(4, 2).Value2 = $row2.B
$ws.Cells.Item(4, 5).Value2 = $row2.E
$ws.Cells.Item(4, 6).Value2 = $row2.F
$ws.Cells.Item(4, 7).Value2 = $row2.G
$ws.Cells.Item(4, 8).Value2 = $row2.H
$ws.Cells.Item(4, 17).Value2 = $row2.Q
$ws.Cells.Item(4, 18).Value2 = $row2.R

# Write old row 4 values into row 5.
$ws.Cells.Item(5, 1).Value2 = $row4.A
$ws.Cells.Item(5, 2).Value2 = $row4.B
$ws.Cells.Item(5, 5).Value2 = $row4.E
$ws.Cells.Item(5, 6).Value2 = $row4.F
$ws.Cells.Item(5, 7).Value2 = $row4.G
$ws.Cells.Item(5, 8).Value2 = $row4.H
$ws.Cells.Item(5, 17).Value2 = $row4.Q
$ws.Cells.Item(5, 18).Value2 = $row4.R
